$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date-like (dd.mm.yyyy) string as literal text, not as a parsed date,
# while restoring the cell's number format to General afterwards so no new style is
# introduced (mirrors existing cells A2:B146 which carry no explicit style override).
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
}

# --- Row 147 : 02.03.2024 (report covering the same figures as 01.03.2024) ---
Set-TextValue 147 1 "02.03.2024"
Set-TextValue 147 2 "01.03.2024"
$ws.Cells.Item(147, 3).Value = 30228
$ws.Cells.Item(147, 4).Value = 12300
$ws.Cells.Item(147, 5).Value = 8400
$ws.Cells.Item(147, 6).Value = 71377
$ws.Cells.Item(147, 7).Value = 8663
$ws.Cells.Item(147, 8).Value = 6327
$ws.Cells.Item(147, 9).Value = 8000
$ws.Cells.Item(147, 10).Value = 417
$ws.Cells.Item(147, 11).Value = 108
$ws.Cells.Item(147, 12).Value = 4600
$ws.Cells.Item(147, 13).Value = "https://web.archive.org/web/20240302133451/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# --- Row 148 : 03.03.2024 ---
Set-TextValue 148 1 "03.03.2024"
Set-TextValue 148 2 "01.03.2024"
$ws.Cells.Item(148, 3).Value = 30228
$ws.Cells.Item(148, 4).Value = 12300
$ws.Cells.Item(148, 5).Value = 8400
$ws.Cells.Item(148, 6).Value = 71377
$ws.Cells.Item(148, 7).Value = 8663
$ws.Cells.Item(148, 8).Value = 6327
$ws.Cells.Item(148, 9).Value = 8000
$ws.Cells.Item(148, 10).Value = 417
$ws.Cells.Item(148, 11).Value = 108
$ws.Cells.Item(148, 12).Value = 4600
$ws.Cells.Item(148, 13).Value = "https://web.archive.org/web/20240303020113/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# --- Row 149 : 04.04.2024 (new tracker figures) ---
Set-TextValue 149 1 "04.04.2024"
Set-TextValue 149 2 "04.04.2024"
$ws.Cells.Item(149, 3).Value = 30534
$ws.Cells.Item(149, 4).Value = 12300
$ws.Cells.Item(149, 5).Value = 8400
$ws.Cells.Item(149, 6).Value = 71920
$ws.Cells.Item(149, 7).Value = 8663
$ws.Cells.Item(149, 8).Value = 6327
$ws.Cells.Item(149, 9).Value = 8000
$ws.Cells.Item(149, 10).Value = 420
$ws.Cells.Item(149, 11).Value = 110
$ws.Cells.Item(149, 12).Value = 4600
$ws.Cells.Item(149, 13).Value = "Updated from the tracker"
